$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet/tab ---
$ws.Name = "Datos Exportados"

# --- Remove the old A1:B1 merge (header now spans 4 distinct columns) ---
$ws.Range("A1:B1").UnMerge()

# --- Row values -----------------------------------------------------
# Header
$ws.Range("A1").Value = "NOMBRE"
$ws.Range("B1").Value = "STOCK ACTUAL"
$ws.Range("C1").Value = "MÁXIMO STOCK"
$ws.Range("D1").Value = "ÚLTIMA ACTUALIZACIÓN"

# Row 2
$ws.Range("A2").Value = "Patio Trasero"
$ws.Range("B2").Value = 535
$ws.Range("C2").Value = 3000
$ws.Range("D2").Value = "21/11/2024"

# Row 3
$ws.Range("A3").Value = "No registrado"
$ws.Range("B3").Value = 824
$ws.Range("C3").Value = 500
$ws.Range("D3").Value = "22/11/2024"

# Row 4
$ws.Range("A4").Value = "asdfsdfs"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = "22/11/2024"

# Row 5
$ws.Range("A5").Value = "asdfsdfs"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "22/11/2024"

# Row 6
$ws.Range("A6").Value = "asdfsdfs"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "22/11/2024"

# Row 7
$ws.Range("A7").Value = "Inventario Secundario"
$ws.Range("B7").Value = 535
$ws.Range("C7").Value = 600
$ws.Range("D7").Value = "22/11/2024"

# Row 8
$ws.Range("A8").Value = "Inventario Principal"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 500
$ws.Range("D8").Value = "22/11/2024"

# --- Formatting -------------------------------------------------------
# Carry the existing banded row styles (s1 header / s2 white / s3 gray)
# across into the two new columns C & D, and onto the brand-new rows
# 6-8, by copy/paste-special(formats) from the matching band in A:B.
$ws.Range("A1:B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$ws.Range("A2:B2").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

$ws.Range("A3:B3").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)

$ws.Range("A4:B4").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

$ws.Range("A3:D3").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

$ws.Range("A4:D4").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)

$ws.Range("A3:D3").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)

$ws.Range("A4:D4").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

# --- Column widths (characters: A=26, B=17, C=17, D=25) ---
$ws.Range("A:A").ColumnWidth = 25.15
$ws.Range("B:C").ColumnWidth = 16.2
$ws.Range("D:D").ColumnWidth = 24.2
